$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.508.55'
$ws.Range('E2').Value = '  +0.93%  '

# Row 3
$ws.Range('D3').Value = '1.626.45'
$ws.Range('E3').Value = '  +1.43%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.56'
$ws.Range('E5').Value = '  +0.15%  '

# Row 6
$ws.Range('E6').Value = '  -0.03%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.487'
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('E8').Value = '  +0.21%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').Value = '  +0.67%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.93'
$ws.Range('E10').Value = '  +4.08%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0828'
$ws.Range('E11').Value = '  +2.12%  '

# Row 12
$ws.Range('D12').Value = '1.851.92'
$ws.Range('E12').Value = '  +1.43%  '

# Row 13
$ws.Range('D13').Value = '1.640.41'
$ws.Range('E13').Value = '  +2.23%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').Value = '  +0.19%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  +1.75%  '

# Row 16
$ws.Range('D16').Value = '26.497.24'
$ws.Range('E16').Value = '  +0.95%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.75'
$ws.Range('E17').Value = '  +2.39%  '

# Row 18
$ws.Range('E18').Value = '  +0.13%  '

# Row 19
$ws.Range('E19').Value = '  -0.11%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.86'
$ws.Range('E20').Value = '  -0.14%  '

# Row 21
$ws.Range('E21').Value = '  -0.02%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.36'
$ws.Range('E22').Value = '  +0.83%  '

# Row 23
$ws.Range('E23').Value = '  +0.74%  '

# Row 24
$ws.Range('E24').Value = '  -3.95%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.12'
$ws.Range('E25').Value = '  +0.33%  '

# Row 26
$ws.Range('E26').Value = '  -0.06%  '

# Row 27
$ws.Range('E27').Value = '  -2.13%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.36'
$ws.Range('E28').Value = '  +1.21%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.62'
$ws.Range('E29').Value = '  +1.14%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0520'
$ws.Range('E30').Value = '  +5.53%  '

# Row 31
$ws.Range('E31').Value = '  +0.42%  '

# Row 32
$ws.Range('E32').Value = '  +1.31%  '

# Row 33
$ws.Range('E33').Value = '  +0.53%  '

# Row 34
$ws.Range('E34').Value = '  +1.33%  '

# Row 35
$ws.Range('E35').Value = '  -0.43%  '

# Row 36
$ws.Range('D36').Value = '1.151.65'
$ws.Range('E36').Value = '  +0.29%  '

# Row 37
$ws.Range('E37').Value = '  +0.75%  '

# Row 38
$ws.Range('E38').Value = '  +2.05%  '

# Row 39
$ws.Range('E39').Value = '  -0.04%  '

# Row 40
$ws.Range('E40').Value = '  -0.38%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.499'
$ws.Range('E41').Value = '  +0.21%  '

# Row 42
$ws.Range('E42').Value = '  +3.48%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.786'
$ws.Range('E43').Value = '  +0.71%  '

# Row 44
$ws.Range('D44').Value = '1.762.53'
$ws.Range('E44').Value = '  +1.31%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.15'
$ws.Range('E45').Value = '  +0.29%  '

# Row 46
$ws.Range('E46').Value = '  +1.13%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.08'
$ws.Range('E47').Value = '  -0.16%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0510'
$ws.Range('E48').Value = '  +0.71%  '

# Row 49
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.409'
$ws.Range('E49').Value = '  +0.60%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.50'
$ws.Range('E50').Value = '  +4.35%  '

# Row 51
$ws.Range('E51').Value = '  -0.15%  '
